{"js": "// The document previously had no word/styles.xml part at all. The edit\n// mints that part by defining the (implicit) \"Normal\" paragraph style\n// explicitly, so the package gains a styles part that declares the\n// default \"Normal\" style - matching the target diff, which adds a new\n// word/styles.xml containing a single paragraph style \"Normal\".\nconst style = context.document.addStyle(\"Normal\", Word.StyleType.paragraph);\nstyle.load(\"nameLocal,type,builtIn\");\nawait context.sync();\n", "ps1": "# The document previously had no word/styles.xml part at all. The edit\n# mints that part by defining the (implicit) \"Normal\" paragraph style\n# explicitly, so the package gains a styles part that declares the\n# default \"Normal\" style - matching the target diff, which adds a new\n# word/styles.xml containing a single paragraph style \"Normal\".\n$d = $word.ActiveDocument\n$style = $d.Styles.Add(\"Normal\", 1)  # wdStyleTypeParagraph = 1\n$style.NameLocal = \"Normal\"\n"}
